# Daily attendance processing - 2025-12-01 17:57:36
#
# The "Recorded By" column (column G) on the "Session Analysis Results"
# sheet lists the users who recorded a given session, separated by ", ".
# For every row whose list does not already start with the literal
# value "System", rotate the list left by one position (the first
# name moves to the end of the list). Rows that already start with
# "System" - and single-value cells - are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $text = $cell.Text

    if ($text -ne $null -and $text -ne "") {
        $parts = $text -split ", "

        if ($parts.Length -gt 1 -and -not $parts[0].Equals("System")) {
            $rotated = $parts[1..($parts.Length - 1)] + $parts[0]
            $cell.Value = ($rotated -join ", ")
        }
    }
}
